# changes in test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Row 2: update name/unit values
$ws.Range("B2").Value = "srusi"
$ws.Range("C2").Value = "ftgi"

# Row 3: clear name/unit, flip selectedField flag to Yes
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("F3").Value = "Yes"

# Row 4: update name/unit values
$ws.Range("B4").Value = "sdit"
$ws.Range("C4").Value = "udyth"

# Update the active selection to match the saved workbook state
$ws.Range("K8").Select()
